$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.154.28'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '1.787.09'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.83'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0688'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '2.044.96'
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.10%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.783.20'
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '34.113.34'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("E32").Value = '  +2.80%  '
$ws.Range("E33").Value = '  +3.54%  '
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").Value = '1.444.62'
$ws.Range("E35").Value = '  +4.74%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.90%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.654'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0191'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.98%  '
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.90%  '
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("E45").Value = '  +4.32%  '
$ws.Range("E46").Value = '  +1.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").Value = '1.946.76'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '105.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("E51").Value = '  +0.03%  '
